# Update gh-pages to output generated at 456a3b4
# This applies numeric updates to the "F" column (想去人数 / interested-count)
# on three worksheets: 展览, 演出, and 全部类型.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 10200
$ws1.Range("F11").Value = 181
$ws1.Range("F23").Value = 1113
$ws1.Range("F26").Value = 624
$ws1.Range("F28").Value = 154
$ws1.Range("F30").Value = 2759
$ws1.Range("F39").Value = 1224
$ws1.Range("F40").Value = 542
$ws1.Range("F41").Value = 5232
$ws1.Range("F46").Value = 36

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 5

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 10200
$ws4.Range("F14").Value = 181
$ws4.Range("F22").Value = 1113
$ws4.Range("F26").Value = 5
$ws4.Range("F27").Value = 624
$ws4.Range("F29").Value = 154
$ws4.Range("F31").Value = 2759
